$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -eq $null) { continue }
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    $count = $parts.Length

    if ($count -gt 1) {
        $reversed = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }
        $newText = [string]::Join(", ", $reversed)
        $cell.Value = $newText
    }
}
